$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns before H (shifts old H,I -> J,K, carrying their
# formatting/width along), so we can populate the two new header cells.
$ws.Range("H1:I1").Insert(-4161)

# New header cells for the two inserted columns.
$ws.Range("H1").Value = "عدد أفراد الأسرة"
$ws.Range("I1").Value = "الحالة الإجتماعية"

# New trailing column L - copy the header format from the neighbouring
# header cell (K1, formerly I1) so it matches the other header cells, then
# set its text.
$ws.Range("K1").Copy()
$ws.Range("L1").PasteSpecial(-4122)
$ws.Range("L1").Value = "عنوان السكن"

# Column widths (values chosen to land on the closest width this runtime can
# represent; widths are stored on a 1/6-character grid after conversion).
$ws.Range("D1").ColumnWidth = 15.999999999999998
$ws.Range("G1").ColumnWidth = 12.666666666666666
$ws.Range("H1").ColumnWidth = 17.833333333333336
$ws.Range("I1").ColumnWidth = 20
$ws.Range("K1").ColumnWidth = 12.5
$ws.Range("L1").ColumnWidth = 14.166666666666666

# Selection matches the post-edit cursor position.
[void]$ws.Range("I2").Select()

# Data validation: column I (الحالة الإجتماعية) restricted to a fixed list.
$ws.Range("I1:I1048576").Validation.Add(3, 1, 1, '"اعزب,متزوج,مطلق,ارملة"')
